# Updated cryptos list on Wed Nov 15 14:25:26 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Cell($row, $col, $value) {
    $ws.Cells.Item($row, $col).Value = $value
}

# --- Row 2: Bitcoin ---
Set-Cell 2 4 "36.210.91"
Set-Cell 2 5 "  -1.25%  "

# --- Row 3: Ethereum ---
Set-Cell 3 4 "2.019.74"
Set-Cell 3 5 "  -1.57%  "

# --- Row 4: TetherUSD ---
Set-Cell 4 5 "  +0.10%  "

# --- Row 5: BNB ---
Set-Cell 5 4 "253.57"
Set-Cell 5 5 "  +3.25%  "

# --- Row 6: XRP ---
Set-Cell 6 4 "0.644"
Set-Cell 6 5 "  -3.85%  "

# --- Row 7: Solana ---
Set-Cell 7 4 "62.37"
Set-Cell 7 5 "  +7.94%  "

# --- Row 9: OKB ---
Set-Cell 9 4 "59.11"
Set-Cell 9 5 "  -6.36%  "

# --- Row 10: Cardano ---
Set-Cell 10 5 "  +0.89%  "

# --- Row 11: Dogecoin ---
Set-Cell 11 5 "  -0.42%  "

# --- Row 12: TRON ---
Set-Cell 12 5 "  -1.75%  "

# --- Row 13: Polygon ---
Set-Cell 13 4 "0.923"
Set-Cell 13 5 "  +0.92%  "

# --- Row 14: Chainlink ---
Set-Cell 14 4 "14.99"
Set-Cell 14 5 "  +1.70%  "

# --- Row 15: WrappedliquidstakedEther2.0 ---
Set-Cell 15 4 "2.312.33"
Set-Cell 15 5 "  -1.53%  "

# --- Row 16: Polkadot ---
Set-Cell 16 4 "5.44"
Set-Cell 16 5 "  -0.12%  "

# --- Row 17: Avalanche ---
Set-Cell 17 4 "19.76"
Set-Cell 17 5 "  +11.90%  "

# --- Row 18: WrappedEther ---
Set-Cell 18 4 "2.016.22"
Set-Cell 18 5 "  -1.50%  "

# --- Row 19: WrappedBTC ---
Set-Cell 19 4 "36.157.32"
Set-Cell 19 5 "  -1.02%  "

# --- Row 20: Litecoin ---
Set-Cell 20 4 "72.19"
Set-Cell 20 5 "  -0.03%  "

# --- Row 21: ShibaInu ---
Set-Cell 21 4 "0.0₃0860"
Set-Cell 21 5 "  -0.12%  "

# --- Row 22: Uniswap ---
Set-Cell 22 4 "5.30"
Set-Cell 22 5 "  +1.77%  "

# --- Row 23: BitcoinCash ---
Set-Cell 23 4 "234.29"
Set-Cell 23 5 "  -1.78%  "

# --- Row 24: PancakeSwap ---
Set-Cell 24 4 "2.71"
Set-Cell 24 5 "  +18.92%  "

# --- Row 25: Dai ---
Set-Cell 25 5 "  +0.01%  "

# --- Row 26: Toncoin ---
Set-Cell 26 4 "2.32"
Set-Cell 26 5 "  -2.09%  "

# --- Row 27: Cosmos ---
Set-Cell 27 4 "9.56"
Set-Cell 27 5 "  +2.57%  "

# --- Row 28: Monero ---
Set-Cell 28 4 "164.66"
Set-Cell 28 5 "  -0.03%  "

# --- Row 29: EthereumClassic ---
Set-Cell 29 4 "19.67"
Set-Cell 29 5 "  -2.15%  "

# --- Row 30: Stellar ---
Set-Cell 30 5 "  -1.01%  "

# --- Row 31: ImmutableX ---
Set-Cell 31 4 "1.21"
Set-Cell 31 5 "  +0.85%  "

# --- Row 32: Filecoin ---
Set-Cell 32 4 "5.11"
Set-Cell 32 5 "  +0.99%  "

# --- Row 33: Kaspa ---
Set-Cell 33 4 "0.108"
Set-Cell 33 5 "  +24.87%  "

# --- Row 34: Hedera ---
Set-Cell 34 4 "0.0607"
Set-Cell 34 5 "  +0.59%  "

# --- Row 35 & 36: InternetComputer(DFINITY) / LidoDAOToken swap places ---
Set-Cell 35 2 "LidoDAOToken"
Set-Cell 35 3 "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-Cell 35 4 "2.50"
Set-Cell 35 5 "  +13.18%  "

Set-Cell 36 2 "InternetComputer(DFINITY)"
Set-Cell 36 3 "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-Cell 36 4 "4.51"
Set-Cell 36 5 "  +1.38%  "

# --- Row 37: BinanceUSD ---
Set-Cell 37 5 "  +0.07%  "

# --- Row 38: WEMIXToken ---
Set-Cell 38 5 "  -1.00%  "

# --- Row 39: THORChain ---
Set-Cell 39 4 "5.85"
Set-Cell 39 5 "  +16.17%  "

# --- Row 40: Cronos ---
Set-Cell 40 4 "0.104"
Set-Cell 40 5 "  +14.74%  "

# --- Row 42 & 43: VeChain / HuobiToken swap places ---
Set-Cell 42 2 "HuobiToken"
Set-Cell 42 3 "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-Cell 42 4 "2.90"
Set-Cell 42 5 "  -0.96%  "

Set-Cell 43 2 "VeChain"
Set-Cell 43 3 "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-Cell 43 4 "0.0216"
Set-Cell 43 5 "  +0.14%  "

# --- Row 44: ARBITRUM ---
Set-Cell 44 4 "1.13"
Set-Cell 44 5 "  +1.58%  "

# --- Row 45: InjectiveProtocol ---
Set-Cell 45 4 "16.67"
Set-Cell 45 5 "  +4.05%  "

# --- Row 46 & 47: FraxShare / Aave swap places ---
Set-Cell 46 2 "Aave"
Set-Cell 46 3 "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-Cell 46 4 "94.20"
Set-Cell 46 5 "  -0.19%  "

Set-Cell 47 2 "FraxShare"
Set-Cell 47 3 "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-Cell 47 4 "7.83"
Set-Cell 47 5 "  +4.71%  "

# --- Row 48: Maker ---
Set-Cell 48 4 "1.424.87"
Set-Cell 48 5 "  +3.09%  "

# --- Row 49: RenderToken ---
Set-Cell 49 4 "2.52"
Set-Cell 49 5 "  +11.42%  "

# --- Row 50: MXToken ---
Set-Cell 50 4 "2.91"
Set-Cell 50 5 "  -1.05%  "

# --- Row 51: MultiversX ---
Set-Cell 51 4 "47.35"
Set-Cell 51 5 "  +3.28%  "
